# Update decorator starter and ender for main flow
# - Refresh Start Time / End Time / Run Time on the "Version" sheet
# - Refresh the BT Idle decorator stats (Average/Max/Min/Sdev/Raw) on the
#   "3_3" and "1_8" sheets to match a re-run with a new starter/ender pair.

$wb = $excel.ActiveWorkbook

# ---- Version sheet: Start/End/Run time -------------------------------
$verSheet = $wb.Worksheets.Item("Version")
$verSheet.Range("C10").Value = "09:11:53 03-22-2017 Wednesday Pacific Daylight Time"
$verSheet.Range("C11").Value = "09:12:28 03-22-2017 Wednesday Pacific Daylight Time"
$verSheet.Range("C12").Value = "0:0:35 0 days"

# ---- "3_3" sheet: BT Idle row (row 3) and row 4 -----------------------
$s33 = $wb.Worksheets.Item("3_3")

$s33.Range("C3").Value = 0.594
$s33.Range("D3").Value = 0.622
$s33.Range("E3").Value = 0.5669999999999999
$s33.Range("F3").Value = 0.014
$s33.Range("H3").Value = "0.601303,0.583211,0.597233,0.616622,0.603171,0.605304,0.581078,0.596391,0.598429,0.600202,0.577077,0.576703,0.602518,0.590073,0.60455,0.601303,0.575786,0.615888,0.603531,0.570392,0.608001,0.598422,0.573082,0.581452,0.60938,0.583015,0.603075,0.610596,0.583945,0.601588,0.609754,0.574292,0.611988,0.595739,0.569937,0.583388,0.606886,0.573925,0.609115,0.616058,0.574115,0.60298,0.608646,0.584686,0.59082,0.600392,0.598055,0.606044,0.595745,0.566955,0.572715,0.607253,0.581717,0.575589,0.605766,0.576805,0.600841,0.610589,0.606696,0.572906,0.603442,0.60762,0.607437,0.575684,0.602518,0.595745,0.574197,0.600936,0.607063,0.574013,0.603904,0.590168,0.584876,0.600385,0.60224,0.56817,0.611887,0.601499,0.582566,0.612084,0.580141,0.607811,0.573463,0.609672,0.607348,0.57709,0.602613,0.580704,0.599631,0.57188,0.598531,0.602233,0.585249,0.621921,0.575494,0.572532,0.598714,0.57902,0.607906,0.605949"

$s33.Range("C4").Value = 9.831
$s33.Range("D4").Value = 9.851000000000001
$s33.Range("E4").Value = 9.81
$s33.Range("F4").Value = 0.008999999999999999
$s33.Range("H4").Value = "9.83645,9.819934,9.825505,9.843508,9.832557,9.821782,9.817516,9.83005,9.831538,9.831721,9.819663,9.820675,9.840723,9.847407,9.827638,9.817142,9.82058,9.836273,9.835798,9.821511,9.818352,9.827733,9.838596,9.832557,9.822713,9.818637,9.838589,9.842869,9.826436,9.821144,9.833386,9.832835,9.835254,9.833487,9.823277,9.836653,9.841096,9.828847,9.819567,9.818542,9.839235,9.840641,9.832285,9.821055,9.836178,9.843888,9.833209,9.82717,9.827638,9.834316,9.842958,9.832835,9.814826,9.830063,9.839425,9.842027,9.828936,9.819928,9.832645,9.851395,9.831721,9.826062,9.818827,9.839989,9.846578,9.833304,9.812883,9.828936,9.836932,9.836545,9.820118,9.822992,9.832183,9.842217,9.831544,9.819656,9.818732,9.840824,9.845559,9.834038,9.822631,9.82522,9.842027,9.840729,9.824751,9.826334,9.831082,9.837564,9.836836,9.816402,9.836368,9.836456,9.842679,9.82774,9.810084,9.830702,9.841375,9.839242,9.830804,9.826436"

# ---- "1_8" sheet: BT Idle row (row 3) and row 4 -----------------------
$s18 = $wb.Worksheets.Item("1_8")

$s18.Range("C3").Value = -0.001
$s18.Range("D3").Value = 0.007
$s18.Range("E3").Value = -0.008
$s18.Range("F3").Value = 0.003
$s18.Range("H3").Value = "-0.000263,-0.004974,-0.002982,-0.000251,0.001692,-0.001071,-0.000799,-0.001089,-0.000781,0.000311,-0.000799,-0.003013,-0.002735,-0.001885,0.000878,-0.003827,-0.004659,-0.005245,-0.001077,-0.005208,-0.005245,0.000317,-0.006904,-0.005208,-0.000768,-0.001607,0.003369,-0.001891,-0.000251,3.9e-05,0.001963,-0.000226,-0.003833,0.000866,0.001421,-0.003278,-0.001077,-0.001046,-0.003266,0.005318,-0.002193,-0.003852,-0.003852,-0.003845,-0.001934,0.002537,-0.002452,-0.001872,-0.00387,-0.00387,0.001149,0.001692,-0.008291,0.006427,0.005848,-0.00358,-0.000497,0.001717,-0.004111,-0.000244,-0.001922,0.000853,-0.000787,-0.003852,-0.004111,0.005614,-0.001897,0.001445,-0.005233,-0.004961,-0.003019,0.001969,-0.000787,-0.002181,-0.002464,-0.003297,0.004448,-0.001064,-0.00022,0.003943,-0.004123,0.001168,-0.008285,0.003955,-0.001897,-0.000238,-0.006065,-0.006645,0.006705,0.00311,-0.004672,-0.002199,-0.001909,-0.001071,0.001131,0.002253,0.000298,-0.003025,0.000582,-0.004968"

$s18.Range("C4").Value = -0.001
$s18.Range("D4").Value = 0.007
$s18.Range("E4").Value = -0.012
$s18.Range("F4").Value = 0.003
$s18.Range("H4").Value = "-0.000214,0.000884,-0.004111,-0.005806,-0.001891,-0.002181,-0.000238,0.000311,-0.005492,-0.002729,0.003116,0.000329,0.003104,-0.001071,-0.005775,0.00393,0.005299,-0.01217,0.002506,-0.003019,0.000304,0.002259,0.003092,0.003943,-0.00136,-0.002458,0.000317,-0.002476,4.5e-05,-0.001342,-0.000793,0.005046,-0.001064,-0.006078,-0.001909,0.002549,0.000872,0.000607,0.00615,-0.000818,-0.002735,0.003369,0.001451,-0.002458,-0.007157,-0.001909,-0.00522,0.003098,-0.003858,-0.007736,-0.000509,-0.005523,0.000329,-0.004129,-0.003845,0.000594,-0.005535,-0.008581,0.001433,0.00282,-0.002458,-0.002994,0.001704,0.003085,0.003665,6.4e-05,0.003918,-0.004154,-0.00136,-0.003038,-0.002187,0.002518,-0.001613,-0.005257,-0.001589,5.8e-05,-0.000226,-0.001077,0.007285,-0.000793,4.5e-05,-0.0044,-0.005227,0.00253,5.2e-05,-0.004419,-0.009407,0.003936,-0.001342,-0.002199,4.5e-05,-0.001064,-0.000509,0.003924,-0.004129,8e-06,0.002228,-0.000226,-0.001354,0.001704"
